$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 166. This pushes the existing rows
# 166-171 ("Chino" garlic entries) down to rows 169-174, matching the
# target layout, and leaves rows 166-168 empty for the new "Rosado"
# garlic entries.
$ws.Rows("166:168").Insert()

# Row 166: Ajo Rosado, 1a (cosecha)
$ws.Range("A166").Value = 9
$ws.Range("B166").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C166").Value = "Metropolitana"
$ws.Range("D166").Value = (Get-Date -Year 2022 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = 100112003
$ws.Range("G166").Value = "Ajo"
$ws.Range("H166").Value = "Rosado"
$ws.Range("I166").Value = "1a (cosecha)"
$ws.Range("J166").Value = 250
$ws.Range("K166").Value = 8000
$ws.Range("L166").Value = 9000
$ws.Range("M166").Value = 8500
$ws.Range("N166").Value = "`$/trenza 50 unidades"
$ws.Range("O166").Value = "Región de O'Higgins"
$ws.Range("P166").Value = 1700
$ws.Range("Q166").Value = 5
$ws.Range("R166").Value = "Hortaliza"

# Row 167: Ajo Rosado, 2a (cosecha)
$ws.Range("A167").Value = 9
$ws.Range("B167").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C167").Value = "Metropolitana"
$ws.Range("D167").Value = (Get-Date -Year 2022 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E167").Value = 13
$ws.Range("F167").Value = 100112003
$ws.Range("G167").Value = "Ajo"
$ws.Range("H167").Value = "Rosado"
$ws.Range("I167").Value = "2a (cosecha)"
$ws.Range("J167").Value = 160
$ws.Range("K167").Value = 6000
$ws.Range("L167").Value = 6500
$ws.Range("M167").Value = 6250
$ws.Range("N167").Value = "`$/trenza 50 unidades"
$ws.Range("O167").Value = "Región de O'Higgins"
$ws.Range("P167").Value = 1250
$ws.Range("Q167").Value = 5
$ws.Range("R167").Value = "Hortaliza"

# Row 168: Ajo Rosado, 3a (cosecha)
$ws.Range("A168").Value = 9
$ws.Range("B168").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C168").Value = "Metropolitana"
$ws.Range("D168").Value = (Get-Date -Year 2022 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = 100112003
$ws.Range("G168").Value = "Ajo"
$ws.Range("H168").Value = "Rosado"
$ws.Range("I168").Value = "3a (cosecha)"
$ws.Range("J168").Value = 97
$ws.Range("K168").Value = 5000
$ws.Range("L168").Value = 5500
$ws.Range("M168").Value = 5247
$ws.Range("N168").Value = "`$/trenza 50 unidades"
$ws.Range("O168").Value = "Región de O'Higgins"
$ws.Range("P168").Value = 1049
$ws.Range("Q168").Value = 5
$ws.Range("R168").Value = "Hortaliza"

# The now-shifted row 169 (previously 166, "1a (guarda)") changes its
# unit-of-sale text and quantity/price columns.
$ws.Range("N169").Value = "`$/trenza 50 unidades"
$ws.Range("Q169").Value = 5
